# Apply the "LinuxForHealth" re-branding edit to the FHIR StructureDefinition workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Metadata" worksheet: update URL, Version, Date and Publisher values.
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/company-code"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---------------------------------------------------------------------------
# 2. "Elements" worksheet: the root "Extension" row (row 2) no longer shows
#    the ele-1 / ext-1 constraint text in the "Constraint(s)" column (AI).
# ---------------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

$elements.Range("AI2").Value = ""
